$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 210, shifting existing rows 210:241 down to 211:242
$ws.Rows.Item(210).Insert()

# Populate the newly inserted row 210 with data (matches the existing pattern of this block)
$ws.Cells.Item(210, 1).Value = 3
$ws.Cells.Item(210, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(210, 3).Value = "Coquimbo"
$ws.Cells.Item(210, 4).Value = 44984
$ws.Cells.Item(210, 4).Style = $ws.Cells.Item(211, 4).Style
$ws.Cells.Item(210, 4).NumberFormat = $ws.Cells.Item(211, 4).NumberFormat
$ws.Cells.Item(210, 5).Value = 5
$ws.Cells.Item(210, 6).Value = 100112052
$ws.Cells.Item(210, 7).Value = "Albahaca"
$ws.Cells.Item(210, 8).Value = "Sin especificar"
$ws.Cells.Item(210, 9).Value = "Primera"
$ws.Cells.Item(210, 10).Value = 60
$ws.Cells.Item(210, 11).Value = 4500
$ws.Cells.Item(210, 12).Value = 4500
$ws.Cells.Item(210, 13).Value = 4500
$ws.Cells.Item(210, 14).Value = "`$/docena de matas"
$ws.Cells.Item(210, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(210, 16).Value = 750
$ws.Cells.Item(210, 17).Value = 6
$ws.Cells.Item(210, 18).Value = "Hortaliza"
